# Auto-generated data refresh for Anima Profits workbook.
# For each affected row, update the H-N (price/profit) columns to the
# refreshed market-data snapshot values. Cells that are newly introduced
# are written directly; cells removed by the refresh are cleared so the
# <c> element disappears entirely (matching upstream export behaviour).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 937.04443
$ws.Range("J17").Value = 1005.6579
$ws.Range("L17").Value = 3016.9737
$ws.Range("N17").Value = -3352.9737
# Row 33
$ws.Range("H33").Value = 1572.2142
$ws.Range("I33").Value = 992.6667
$ws.Range("J33").Value = 5049.5
$ws.Range("K33").Value = 992.6667
$ws.Range("L33").Value = 5049.5
$ws.Range("M33").Value = -763.6667
$ws.Range("N33").Value = -5507.5
# Row 41
$ws.Range("H41").Value = 180
$ws.Range("I41").Value = 200
$ws.Range("K41").Value = 200
$ws.Range("M41").Value = 240
# Row 62
$ws.Range("H62").Value = 3350
$ws.Range("I62").Value = 3133.3333
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3133.3333
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -2509.3333
$ws.Range("N62").Value = -5248
# Row 65
$ws.Range("H65").Value = 3350
$ws.Range("I65").Value = 3133.3333
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 15666.6665
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -12546.6665
$ws.Range("N65").Value = -26240
# Row 121
$ws.Range("H121").Value = 4800
$ws.Range("J121").Value = 4800
$ws.Range("L121").Value = 14400
$ws.Range("N121").Value = -17894
# Row 129
$ws.Range("H129").Value = 1160.5428
$ws.Range("I129").Value = 463.875
$ws.Range("J129").Value = 1747.2106
$ws.Range("K129").Value = 1391.625
$ws.Range("L129").Value = 5241.6318
$ws.Range("M129").Value = 3608.375
$ws.Range("N129").Value = -15241.6318
# Row 131
$ws.Range("H131").Value = 1491.1818
$ws.Range("I131").Value = 287.5
$ws.Range("K131").Value = 862.5
$ws.Range("M131").Value = 4177.5
# Row 138
$ws.Range("H138").Value = 1978.0714
$ws.Range("I138").Value = 1509.3
$ws.Range("J138").Value = 2518.9614
$ws.Range("K138").Value = 4527.9
$ws.Range("L138").Value = 7556.8842
$ws.Range("M138").Value = 612.1000000000004
$ws.Range("N138").Value = -17836.8842

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 574201.5600000001
$ws.Range("I32").Value = 654043.4
$ws.Range("J32").Value = 25289
$ws.Range("K32").Value = 654043.4
$ws.Range("L32").Value = 25289
$ws.Range("M32").Value = -653756.4
$ws.Range("N32").Value = -25863
# Row 34
$ws.Range("H34").Value = 53771
$ws.Range("I34").Value = 5000
$ws.Range("K34").Value = 5000
$ws.Range("M34").Value = -4729
# Row 61
$ws.Range("H61").Value = 2094.0417
$ws.Range("I61").Value = 2094.0417
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2094.0417
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1882.0417
$ws.Range("N61").ClearContents()
# Row 63
$ws.Range("H63").Value = 5682.091
$ws.Range("I63").Value = 5002
$ws.Range("J63").Value = 6248.8335
$ws.Range("K63").Value = 5002
$ws.Range("L63").Value = 6248.8335
$ws.Range("M63").Value = -4316
$ws.Range("N63").Value = -7620.8335
# Row 66
$ws.Range("H66").Value = 5682.091
$ws.Range("I66").Value = 5002
$ws.Range("J66").Value = 6248.8335
$ws.Range("K66").Value = 25010
$ws.Range("L66").Value = 31244.1675
$ws.Range("M66").Value = -21578
$ws.Range("N66").Value = -38108.1675
# Row 88
$ws.Range("H88").Value = 2661.5
$ws.Range("J88").Value = 2775
$ws.Range("L88").Value = 2775
$ws.Range("N88").Value = -3587
# Row 91
$ws.Range("H91").Value = 2661.5
$ws.Range("J91").Value = 2775
$ws.Range("L91").Value = 2775
$ws.Range("N91").Value = -5583
# Row 97
$ws.Range("H97").Value = 791.25
$ws.Range("I97").Value = 791.25
$ws.Range("K97").Value = 791.25
$ws.Range("M97").Value = -295.25
# Row 102
$ws.Range("H102").Value = 3010
$ws.Range("I102").Value = 3111.111
$ws.Range("K102").Value = 3111.111
$ws.Range("M102").Value = -1489.111
# Row 136
$ws.Range("H136").Value = 2094.0417
$ws.Range("I136").Value = 2094.0417
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6282.125100000001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3732.125100000001
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1996.8334
$ws.Range("I86").Value = 1911.625
$ws.Range("J86").Value = 2167.25
$ws.Range("K86").Value = 1911.625
$ws.Range("L86").Value = 2167.25
$ws.Range("M86").Value = -788.625
$ws.Range("N86").Value = -4413.25
# Row 89
$ws.Range("H89").Value = 1996.8334
$ws.Range("I89").Value = 1911.625
$ws.Range("J89").Value = 2167.25
$ws.Range("K89").Value = 9558.125
$ws.Range("L89").Value = 10836.25
$ws.Range("M89").Value = -3942.125
$ws.Range("N89").Value = -22068.25

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 700
$ws.Range("I16").Value = 685.7143
$ws.Range("J16").Value = 750
$ws.Range("K16").Value = 685.7143
$ws.Range("L16").Value = 750
$ws.Range("M16").Value = -398.7143
$ws.Range("N16").Value = -1324
# Row 31
$ws.Range("H31").Value = 4996.9253
$ws.Range("I31").Value = 1177.027
$ws.Range("J31").Value = 9708.134
$ws.Range("K31").Value = 1177.027
$ws.Range("L31").Value = 9708.134
$ws.Range("M31").Value = -882.027
$ws.Range("N31").Value = -10298.134
# Row 33
$ws.Range("H33").Value = 1000
$ws.Range("I33").Value = 1000
$ws.Range("K33").Value = 1000
$ws.Range("M33").Value = -621
# Row 34
$ws.Range("H34").Value = 4996.9253
$ws.Range("I34").Value = 1177.027
$ws.Range("J34").Value = 9708.134
$ws.Range("K34").Value = 1177.027
$ws.Range("L34").Value = 9708.134
$ws.Range("M34").Value = -975.027
$ws.Range("N34").Value = -10112.134
# Row 58
$ws.Range("H58").Value = 1357.15
$ws.Range("I58").Value = 1240.5
$ws.Range("J58").Value = 2407
$ws.Range("K58").Value = 1240.5
$ws.Range("L58").Value = 2407
$ws.Range("M58").Value = -1037.5
$ws.Range("N58").Value = -2813
# Row 107
$ws.Range("H107").Value = 5209331.5
$ws.Range("I107").Value = 6945220
$ws.Range("K107").Value = 6945220
$ws.Range("M107").Value = -6943300
# Row 113
$ws.Range("H113").Value = 700
$ws.Range("I113").Value = 685.7143
$ws.Range("J113").Value = 750
$ws.Range("K113").Value = 685.7143
$ws.Range("L113").Value = 750
$ws.Range("M113").Value = 1484.2857
$ws.Range("N113").Value = -5090
# Row 136
$ws.Range("H136").Value = 1357.15
$ws.Range("I136").Value = 1240.5
$ws.Range("J136").Value = 2407
$ws.Range("K136").Value = 3721.5
$ws.Range("L136").Value = 7221
$ws.Range("M136").Value = -1171.5
$ws.Range("N136").Value = -12321

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 7600
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
# Row 59
$ws.Range("H59").Value = 2824.95
$ws.Range("I59").Value = 1833
$ws.Range("K59").Value = 5499
$ws.Range("M59").Value = -4959
# Row 131
$ws.Range("H131").Value = 1043.3684
$ws.Range("J131").Value = 1134.7878
$ws.Range("L131").Value = 3404.3634
$ws.Range("N131").Value = -13484.3634
# Row 137
$ws.Range("H137").Value = 11123352
$ws.Range("J137").Value = 5717.75
$ws.Range("L137").Value = 17153.25
$ws.Range("N137").Value = -27353.25

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1618.8
$ws.Range("I97").Value = 1333.5294
$ws.Range("K97").Value = 1333.5294
$ws.Range("M97").Value = -837.5293999999999
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2940.4
$ws.Range("I7").Value = 2940.4
$ws.Range("K7").Value = 2940.4
$ws.Range("M7").Value = -2828.4
# Row 22
$ws.Range("H22").Value = 6158.125
$ws.Range("I22").Value = 1444.4445
$ws.Range("K22").Value = 1444.4445
$ws.Range("M22").Value = -1149.4445
# Row 27
$ws.Range("H27").Value = 6158.125
$ws.Range("I27").Value = 1444.4445
$ws.Range("K27").Value = 1444.4445
$ws.Range("M27").Value = -1337.4445
# Row 31
$ws.Range("H31").Value = 1361.5
$ws.Range("I31").Value = 1282
$ws.Range("J31").Value = 1600
$ws.Range("K31").Value = 1282
$ws.Range("L31").Value = 1600
$ws.Range("M31").Value = -1034
$ws.Range("N31").Value = -2096
# Row 55
$ws.Range("H55").Value = 842.2941
$ws.Range("I55").Value = 667.375
$ws.Range("K55").Value = 667.375
$ws.Range("M55").Value = -494.375
# Row 61
$ws.Range("H61").Value = 4560
$ws.Range("J61").Value = 5300
$ws.Range("L61").Value = 5300
$ws.Range("N61").Value = -5704
# Row 113
$ws.Range("H113").Value = 4560
$ws.Range("J113").Value = 5300
$ws.Range("L113").Value = 5300
$ws.Range("N113").Value = -9640
# Row 123
$ws.Range("H123").Value = 28800
$ws.Range("J123").Value = 28800
$ws.Range("L123").Value = 28800
$ws.Range("N123").Value = -38600
# Row 126
$ws.Range("H126").Value = 2940.4
$ws.Range("I126").Value = 2940.4
$ws.Range("K126").Value = 8821.200000000001
$ws.Range("M126").Value = -6351.200000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 76
$ws.Range("H76").Value = 68333.336
$ws.Range("J76").Value = 68333.336
$ws.Range("L76").Value = 68333.336
$ws.Range("N76").Value = -68963.336
# Row 79
$ws.Range("H79").Value = 68333.336
$ws.Range("J79").Value = 68333.336
$ws.Range("L79").Value = 68333.336
$ws.Range("N79").Value = -70517.336
# Row 132
$ws.Range("H132").Value = 6175015
$ws.Range("I132").Value = 2631.8
$ws.Range("J132").Value = 9805829
$ws.Range("K132").Value = 7895.400000000001
$ws.Range("L132").Value = 29417487
$ws.Range("M132").Value = -5365.400000000001
$ws.Range("N132").Value = -29422547
# Row 136
$ws.Range("H136").Value = 1969.7073
$ws.Range("I136").Value = 1524.5555
$ws.Range("J136").Value = 2828.2144
$ws.Range("K136").Value = 4573.666499999999
$ws.Range("L136").Value = 8484.643199999999
$ws.Range("M136").Value = -2023.666499999999
$ws.Range("N136").Value = -13584.6432
